$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.960561752319336
$ws.Range("B1").Value = 1.947490215301514
$ws.Range("C1").Value = 2.362796306610107
$ws.Range("D1").Value = 2.36708402633667
$ws.Range("E1").Value = 1.815576672554016
